$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Sheet "Restricciones_del_lider"
$ws2 = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $ws2.Range("A2") "1.0499999999999998 - x"
Set-TextValue $ws2.Range("B2") "-2.05"
Set-TextValue $ws2.Range("D2") "0.24"
Set-TextValue $ws2.Range("A3") "-1.05 + x"
Set-TextValue $ws2.Range("B3") "0.050000000000000044"
Set-TextValue $ws2.Range("D3") "0.72"

# Sheet "Restricciones_del_follower"
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws3.Range("D2") "0.47"
Set-TextValue $ws3.Range("E2") "3.1"
Set-TextValue $ws3.Range("A3") "0"
Set-TextValue $ws3.Range("B3") "-1"
$ws3.Range("C3").Value = "J_0_LP_v"
Set-TextValue $ws3.Range("D3") "0.88"
Set-TextValue $ws3.Range("E3") "0"
Set-TextValue $ws3.Range("F3") "6.0"

# Sheet "Punto_modificado"
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "1.05"
Set-TextValue $ws4.Range("B2") "2.85"

# Sheet "Vector_bf" (index 5) -- name collides case-insensitively with "Vector_BF",
# so address both of these two sheets by their 1-based index instead of by name.
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-1.834875"

# Sheet "Vector_BF" (index 6)
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-0.12122500000000014"
Set-TextValue $ws6.Range("A3") "-2.1185"

# Sheet "Vector_Alpha" - A2 stays a true numeric cell
$ws7 = $wb.Worksheets.Item("Vector_Alpha")
$ws7.Range("A2").Value = 2.9699999999999998
